$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header values in B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2 and D2 are cleared (deleted), C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -1.1348007166966199
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -1.2061323573086691

# Row 3: update B3:E3 values
$ws.Range("B3").Value = -1.5116290300329904
$ws.Range("C3").Value = -0.38332645618870897
$ws.Range("D3").Value = -2.1901451881043488
$ws.Range("E3").Value = 1.8550160902001722

# Update the selection to match the edited range
$ws.Range("B1:E3").Select()
